$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page one")

# Swap the two objective rows so that "More 2 or more people with last
# name beginning with C" now appears before "Birth to last vote average
# less than 100" (both keep their "No" result in column B).
$ws.Range("A3").Value = "More 2 or more people with last name beginning with C"
$ws.Range("A4").Value = "Birth to last vote average less than 100"
